$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells in row 1 (F1:H1), matching the style of the existing headers ---
$ws.Range("F1").Value = "id"
$ws.Range("G1").Value = "source_file"
$ws.Range("H1").Value = "text"

# Copy the header formatting (bold, border, centered) from an existing header cell
$ws.Range("A1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- New data row 2 ---
$ws.Range("A2").Value = "Zhuanlan Sun"

# B2 holds the text "1" (not a number) in the source data; force text storage
# then drop back to the default (unstyled) cell format
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "1"
$ws.Range("B2").Style = "Normal"

$ws.Range("C2").Value = "无"
$ws.Range("D2").Value = "SUG"
$ws.Range("E2").Value = "WRI"
$ws.Range("F2").Value = "c8048836-24fe-4e27-95aa-c7cfb58ac155"
$ws.Range("G2").Value = "rkc_hGb0Z_annotated.xlsx"
$ws.Range("H2").Value = "The structure of the global policies used in the experiments should be mentioned somewhere."
